# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (Especial/Primera/Segunda) for
# Feria Lagunitas de Puerto Montt - Palta, dated 45013 (2023-03-28),
# pushing the existing rows 720-789 down to 723-792.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 720, shifting everything
# below (including the dimension) down by three rows.
$ws.Range("A720:T722").EntireRow.Insert()

# Row 720: Especial
$ws.Cells.Item(720, 1).Value = 4
$ws.Cells.Item(720, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(720, 3).Value = "Los Lagos"
$ws.Cells.Item(720, 4).Value = 45013
$ws.Cells.Item(720, 5).Value = 10
$ws.Cells.Item(720, 6).Value = "Fruta"
$ws.Cells.Item(720, 7).Value = 100106
$ws.Cells.Item(720, 8).Value = "Oleaginosos"
$ws.Cells.Item(720, 9).Value = 100106002
$ws.Cells.Item(720, 10).Value = "Palta"
$ws.Cells.Item(720, 11).Value = "Hass"
$ws.Cells.Item(720, 12).Value = "Especial"
$ws.Cells.Item(720, 13).Value = 200
$ws.Cells.Item(720, 14).Value = 6000
$ws.Cells.Item(720, 15).Value = 6000
$ws.Cells.Item(720, 16).Value = 6000
$ws.Cells.Item(720, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(720, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(720, 19).Value = 6000
$ws.Cells.Item(720, 20).Value = 1

# Row 721: Primera
$ws.Cells.Item(721, 1).Value = 4
$ws.Cells.Item(721, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(721, 3).Value = "Los Lagos"
$ws.Cells.Item(721, 4).Value = 45013
$ws.Cells.Item(721, 5).Value = 10
$ws.Cells.Item(721, 6).Value = "Fruta"
$ws.Cells.Item(721, 7).Value = 100106
$ws.Cells.Item(721, 8).Value = "Oleaginosos"
$ws.Cells.Item(721, 9).Value = 100106002
$ws.Cells.Item(721, 10).Value = "Palta"
$ws.Cells.Item(721, 11).Value = "Hass"
$ws.Cells.Item(721, 12).Value = "Primera"
$ws.Cells.Item(721, 13).Value = 200
$ws.Cells.Item(721, 14).Value = 5700
$ws.Cells.Item(721, 15).Value = 5700
$ws.Cells.Item(721, 16).Value = 5700
$ws.Cells.Item(721, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(721, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(721, 19).Value = 5700
$ws.Cells.Item(721, 20).Value = 1

# Row 722: Segunda
$ws.Cells.Item(722, 1).Value = 4
$ws.Cells.Item(722, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(722, 3).Value = "Los Lagos"
$ws.Cells.Item(722, 4).Value = 45013
$ws.Cells.Item(722, 5).Value = 10
$ws.Cells.Item(722, 6).Value = "Fruta"
$ws.Cells.Item(722, 7).Value = 100106
$ws.Cells.Item(722, 8).Value = "Oleaginosos"
$ws.Cells.Item(722, 9).Value = 100106002
$ws.Cells.Item(722, 10).Value = "Palta"
$ws.Cells.Item(722, 11).Value = "Hass"
$ws.Cells.Item(722, 12).Value = "Segunda"
$ws.Cells.Item(722, 13).Value = 200
$ws.Cells.Item(722, 14).Value = 5500
$ws.Cells.Item(722, 15).Value = 5500
$ws.Cells.Item(722, 16).Value = 5500
$ws.Cells.Item(722, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(722, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(722, 19).Value = 5500
$ws.Cells.Item(722, 20).Value = 1
